$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.489904854225173
$ws.Range("B2").Value = -1.30879001124827

$ws.Range("A3").Value = -0.484597602834602
$ws.Range("B3").Value = -0.559877243650435

$ws.Range("A4").Value = -0.8461102430395187
$ws.Range("B4").Value = -0.6769899601950435

$ws.Range("A5").Value = -0.7454960035135753
$ws.Range("B5").Value = -0.646893511556516

$ws.Range("A6").Value = 0.8213203197013023
$ws.Range("B6").Value = 0.6030301903909269

$ws.Range("A7").Value = -0.08440379268880797
$ws.Range("B7").Value = 0.0460360109313537

$ws.Range("A8").Value = 0.8022037793066553
$ws.Range("B8").Value = 0.53485774027454
